# Scheduled-runner market data refresh.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) on each job sheet with freshly pulled Market Board values.
$wb = $excel.ActiveWorkbook

# --- ALC sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# row 33: H33,I33,J33,K33,L33,M33,N33
$ws.Range("H33").Value = 37.81818
$ws.Range("I33").Value = 40.666668
$ws.Range("J33").Value = 25
$ws.Range("K33").Value = 40.666668
$ws.Range("L33").Value = 25
$ws.Range("M33").Value = 188.333332
$ws.Range("N33").Value = -483
# row 134: H134,J134,L134,N134
$ws.Range("H134").Value = 53081.816
$ws.Range("J134").Value = 53081.816
$ws.Range("L134").Value = 53081.816
$ws.Range("N134").Value = -63221.816
# row 135: H135,I135,J135,K135,L135,M135,N135
$ws.Range("H135").Value = 1733.2632
$ws.Range("I135").Value = 1742.5714
$ws.Range("J135").Value = 1707.2
$ws.Range("K135").Value = 15683.1426
$ws.Range("L135").Value = 15364.8
$ws.Range("M135").Value = -13148.1426
$ws.Range("N135").Value = -20434.8
# row 138: H138,I138,J138,K138,L138,M138,N138
$ws.Range("H138").Value = 4117311.8
$ws.Range("I138").Value = 1195256.1
$ws.Range("J138").Value = 5652629
$ws.Range("K138").Value = 3585768.3
$ws.Range("L138").Value = 16957887
$ws.Range("M138").Value = -3580628.3
$ws.Range("N138").Value = -16968167

# --- ARM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 3214.2104
$ws.Range("I61").Value = 2409.2
$ws.Range("J61").Value = 4108.6665
$ws.Range("K61").Value = 2409.2
$ws.Range("L61").Value = 4108.6665
$ws.Range("M61").Value = -2197.2
$ws.Range("N61").Value = -4532.6665
# row 74: H74,I74,K74,M74
$ws.Range("H74").Value = 17663.215
$ws.Range("I74").Value = 2344.0908
$ws.Range("K74").Value = 2344.0908
$ws.Range("M74").Value = -1470.0908
# row 77: H77,I77,K77,M77
$ws.Range("H77").Value = 17663.215
$ws.Range("I77").Value = 2344.0908
$ws.Range("K77").Value = 11720.454
$ws.Range("M77").Value = -7352.454
# row 110: H110,I110,J110,K110,L110,M110,N110
$ws.Range("H110").Value = 790.0909
$ws.Range("I110").Value = 810.1111
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 810.1111
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 1234.8889
$ws.Range("N110").Value = -4790
# row 122: H122,I122,K122,M122
$ws.Range("H122").Value = 1710
$ws.Range("I122").Value = 1420.1818
$ws.Range("K122").Value = 4260.5454
$ws.Range("M122").Value = -1810.5454
# row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 3214.2104
$ws.Range("I136").Value = 2409.2
$ws.Range("J136").Value = 4108.6665
$ws.Range("K136").Value = 7227.599999999999
$ws.Range("L136").Value = 12325.9995
$ws.Range("M136").Value = -4677.599999999999
$ws.Range("N136").Value = -17425.9995
# row 139: H139,J139,L139,N139
$ws.Range("H139").Value = 42965.8
$ws.Range("J139").Value = 42965.8
$ws.Range("L139").Value = 42965.8
$ws.Range("N139").Value = -53245.8

# --- BSM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# row 20: H20,I20,J20,K20,L20,M20,N20
$ws.Range("H20").Value = 2309.2856
$ws.Range("I20").Value = 2616.25
$ws.Range("J20").Value = 1900
$ws.Range("K20").Value = 2616.25
$ws.Range("L20").Value = 1900
$ws.Range("M20").Value = -2369.25
$ws.Range("N20").Value = -2394
# row 81: H81,J81,L81,N81
$ws.Range("H81").Value = 34573.332
$ws.Range("J81").Value = 34573.332
$ws.Range("L81").Value = 34573.332
$ws.Range("N81").Value = -36695.332
# row 84: H84,J84,L84,N84
$ws.Range("H84").Value = 34573.332
$ws.Range("J84").Value = 34573.332
$ws.Range("L84").Value = 103719.996
$ws.Range("N84").Value = -114327.996
# row 94: H94,I94,K94,M94
$ws.Range("H94").Value = 966.8570999999999
$ws.Range("I94").Value = 984.4737
$ws.Range("K94").Value = 984.4737
$ws.Range("M94").Value = -533.4737

# --- CRP sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# row 3: H3,I3,K3,M3
$ws.Range("H3").Value = 3346750.5
$ws.Range("I3").Value = 3346750.5
$ws.Range("K3").Value = 3346750.5
$ws.Range("M3").Value = -3346637.5
# row 20: H20,J20,L20,N20
$ws.Range("H20").Value = 56555.145
$ws.Range("J20").Value = 56555.145
$ws.Range("L20").Value = 56555.145
$ws.Range("N20").Value = -57027.145
# row 30: H30,J30,L30,N30
$ws.Range("H30").Value = 56555.145
$ws.Range("J30").Value = 56555.145
$ws.Range("L30").Value = 56555.145
$ws.Range("N30").Value = -56737.145
# row 58: H58,I58,J58,K58,L58,M58,N58
$ws.Range("H58").Value = 3023.6924
$ws.Range("I58").Value = 2086.75
$ws.Range("J58").Value = 3440.111
$ws.Range("K58").Value = 2086.75
$ws.Range("L58").Value = 3440.111
$ws.Range("M58").Value = -1883.75
$ws.Range("N58").Value = -3846.111
# row 99: H99,I99,J99,K99,L99,M99,N99
$ws.Range("H99").Value = 8929843
$ws.Range("I99").Value = 20834400
$ws.Range("J99").Value = 1425
$ws.Range("K99").Value = 20834400
$ws.Range("L99").Value = 1425
$ws.Range("M99").Value = -20832902
$ws.Range("N99").Value = -4421
# row 106: H106,J106,L106,N106
$ws.Range("H106").Value = 28500
$ws.Range("J106").Value = 28500
$ws.Range("L106").Value = 28500
$ws.Range("N106").Value = -31024
# row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 8929843
$ws.Range("I126").Value = 20834400
$ws.Range("J126").Value = 1425
$ws.Range("K126").Value = 62503200
$ws.Range("L126").Value = 4275
$ws.Range("M126").Value = -62500730
$ws.Range("N126").Value = -9215
# row 128: H128,J128,L128,N128
$ws.Range("H128").Value = 56555.145
$ws.Range("J128").Value = 56555.145
$ws.Range("L128").Value = 56555.145
$ws.Range("N128").Value = -66515.14499999999
# row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 3743.842
$ws.Range("I132").Value = 3117
$ws.Range("J132").Value = 4818.4287
$ws.Range("K132").Value = 9351
$ws.Range("L132").Value = 14455.2861
$ws.Range("M132").Value = -6821
$ws.Range("N132").Value = -19515.2861
# row 134: H134,I134,J134,K134,L134,M134,N134
$ws.Range("H134").Value = 3606.5217
$ws.Range("I134").Value = 1850.5714
$ws.Range("J134").Value = 6338
$ws.Range("K134").Value = 5551.7142
$ws.Range("L134").Value = 19014
$ws.Range("M134").Value = -3016.7142
$ws.Range("N134").Value = -24084
# row 135: H135,J135,L135,N135
$ws.Range("H135").Value = 42544
$ws.Range("J135").Value = 42544
$ws.Range("L135").Value = 42544
$ws.Range("N135").Value = -52684
# row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 3023.6924
$ws.Range("I136").Value = 2086.75
$ws.Range("J136").Value = 3440.111
$ws.Range("K136").Value = 6260.25
$ws.Range("L136").Value = 10320.333
$ws.Range("M136").Value = -3710.25
$ws.Range("N136").Value = -15420.333
# row 137: H137,I137,J137,K137,L137,M137,N137
$ws.Range("H137").Value = 35837.8
$ws.Range("I137").Value = 9854.5
$ws.Range("J137").Value = 53160
$ws.Range("K137").Value = 9854.5
$ws.Range("L137").Value = 53160
$ws.Range("M137").Value = -4754.5
$ws.Range("N137").Value = -63360
# row 138: H138,J138,L138,N138
$ws.Range("H138").Value = 51000
$ws.Range("J138").Value = 51000
$ws.Range("L138").Value = 51000
$ws.Range("N138").Value = -61280
# row 140: H140,J140,L140,N140
$ws.Range("H140").Value = 57024
$ws.Range("J140").Value = 57024
$ws.Range("L140").Value = 57024
$ws.Range("N140").Value = -67384

# --- CUL sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# row 22: H22,I22,J22,K22,L22,M22,N22
$ws.Range("H22").Value = 1275.25
$ws.Range("I22").Value = 680
$ws.Range("J22").Value = 2267.3333
$ws.Range("K22").Value = 2040
$ws.Range("L22").Value = 6801.999899999999
$ws.Range("M22").Value = -1871
$ws.Range("N22").Value = -7139.999899999999
# row 27: H27,I27,J27,K27,L27,M27,N27
$ws.Range("H27").Value = 1275.25
$ws.Range("I27").Value = 680
$ws.Range("J27").Value = 2267.3333
$ws.Range("K27").Value = 2040
$ws.Range("L27").Value = 6801.999899999999
$ws.Range("M27").Value = -1938
$ws.Range("N27").Value = -7005.999899999999
# row 81: H81,I81,J81,K81,L81,M81,N81
$ws.Range("H81").Value = 7395
$ws.Range("I81").Value = 1806.25
$ws.Range("J81").Value = 29750
$ws.Range("K81").Value = 5418.75
$ws.Range("L81").Value = 89250
$ws.Range("M81").Value = -4295.75
$ws.Range("N81").Value = -91496
# row 84: H84,I84,J84,K84,L84,M84,N84
$ws.Range("H84").Value = 7395
$ws.Range("I84").Value = 1806.25
$ws.Range("J84").Value = 29750
$ws.Range("K84").Value = 16256.25
$ws.Range("L84").Value = 267750
$ws.Range("M84").Value = -10640.25
$ws.Range("N84").Value = -278982
# row 97: H97,J97,L97,N97
$ws.Range("H97").Value = 1110.9166
$ws.Range("J97").Value = 992.5714
$ws.Range("L97").Value = 2977.7142
$ws.Range("N97").Value = -3969.7142
# row 129: H129,I129,J129,K129,L129,M129,N129
$ws.Range("H129").Value = 1161.5385
$ws.Range("I129").Value = 407.14285
$ws.Range("J129").Value = 2041.6666
$ws.Range("K129").Value = 1221.42855
$ws.Range("L129").Value = 6124.9998
$ws.Range("M129").Value = 3778.57145
$ws.Range("N129").Value = -16124.9998
# row 131: H131,I131,J131,K131,L131,M131,N131
$ws.Range("H131").Value = 1460.1666
$ws.Range("I131").Value = 348.25
$ws.Range("J131").Value = 1631.2307
$ws.Range("K131").Value = 1044.75
$ws.Range("L131").Value = 4893.6921
$ws.Range("M131").Value = 3995.25
$ws.Range("N131").Value = -14973.6921

# --- GSM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# row 113: H113,I113,K113,M113
$ws.Range("H113").Value = 1221.8889
$ws.Range("I113").Value = 1112.5264
$ws.Range("K113").Value = 1112.5264
$ws.Range("M113").Value = 1057.4736
# row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 4124.919
$ws.Range("I132").Value = 4012.6155
$ws.Range("K132").Value = 12037.8465
$ws.Range("M132").Value = -9507.8465
# row 133: H133,J133,L133,N133
$ws.Range("H133").Value = 15654.444
$ws.Range("J133").Value = 15654.444
$ws.Range("L133").Value = 15654.444
$ws.Range("N133").Value = -25774.444

# --- LTW sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# row 7: H7,I7,J7,K7,L7,M7,N7
$ws.Range("H7").Value = 3156.25
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 3541.6667
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 3541.6667
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -3765.6667
# row 22: H22,I22,J22,K22,L22,M22,N22
$ws.Range("H22").Value = 1158.3334
$ws.Range("I22").Value = 100
$ws.Range("J22").Value = 1370
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 1370
$ws.Range("M22").Value = 195
$ws.Range("N22").Value = -1960
# row 27: H27,I27,J27,K27,L27,M27,N27
$ws.Range("H27").Value = 1158.3334
$ws.Range("I27").Value = 100
$ws.Range("J27").Value = 1370
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = 1370
$ws.Range("M27").Value = 7
$ws.Range("N27").Value = -1584
# row 36: H36,J36,L36,N36
$ws.Range("H36").Value = 27240
$ws.Range("J36").Value = 27240
$ws.Range("L36").Value = 27240
$ws.Range("N36").Value = -28364
# row 40: H40,I40,K40,M40
$ws.Range("H40").Value = 3081.9092
$ws.Range("I40").Value = 1660.4
$ws.Range("K40").Value = 1660.4
$ws.Range("M40").Value = -1524.4
# row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 3156.25
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3541.6667
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 10625.0001
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -15565.0001

# --- WVR sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# row 101: H101,J101,L101,N101
$ws.Range("H101").Value = 22899.666
$ws.Range("J101").Value = 22899.666
$ws.Range("L101").Value = 22899.666
$ws.Range("N101").Value = -29389.666
# row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 35894.766
$ws.Range("I122").Value = 49510
$ws.Range("J122").Value = 4125.8887
$ws.Range("K122").Value = 148530
$ws.Range("L122").Value = 12377.6661
$ws.Range("M122").Value = -146080
$ws.Range("N122").Value = -17277.6661
